$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new price-history row right after the last used row (row 37 ->
# new row 38), same as the daily scraper does.
$row = $ws.UsedRange.Rows.Count + 1

# Force these new cells to be written as text (matching the existing
# column data, which is stored as shared strings even for numeric-looking
# values), then restore the default "Normal" style so no stray per-cell
# formatting is left behind.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 4).NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2026-02-07"
$ws.Cells.Item($row, 2).Value = "185600"
$ws.Cells.Item($row, 3).Value = "6"
$ws.Cells.Item($row, 4).Value = "0"

$ws.Cells.Item($row, 1).Style = "Normal"
$ws.Cells.Item($row, 2).Style = "Normal"
$ws.Cells.Item($row, 3).Style = "Normal"
$ws.Cells.Item($row, 4).Style = "Normal"
